$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (defer D1 until after the data rows so the shared-string
# table ends up ordered the same way the authoring tool produced it)
$ws.Range("A1").Value = "Test #"
$ws.Range("B1").Value = "MIDI Note"
$ws.Range("C1").Value = "Velocity"
$ws.Range("E1").Value = "Comments / Observations"

# Update existing data row (was "ALL" / "around 1 micro second..."), now first data row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "C4"
$ws.Range("C2").Value = 90
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "snaping back and forth between two "

# Now set the latency header (new shared string appended after C4/snaping)
$ws.Range("D1").Value = "Measured Latency (micro seconds)"

# Add new second data row
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "C4/D4"
$ws.Range("C3").Value = 90
$ws.Range("D3").Value = 2
